# Apply the commit: add a "Table_2" worksheet with capital adequacy ratio
# data, and clear out a few stray empty placeholder cells on "Table_1".

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Clean up empty inline-string placeholder cells on Table_1 ---
$ws1.Range("B2").ClearContents()
$ws1.Range("A3").ClearContents()
$ws1.Range("B37").ClearContents()

# --- Add the new Table_2 worksheet right after Table_1 ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Table_2"

# Row labels (plain text, never mistaken for numbers)
$ws2.Range("A1").Value = "Əmsal"
$ws2.Range("A2").Value = "9.  I dərəcəli  kapitalın  adekvatlıq əmsalı"
$ws2.Range("A3").Value = "10. məcmu kapitalın  adekvatlıq  əmsalı"
$ws2.Range("A4").Value = "11. Leverec əmsalı"

$ws2.Range("B1").Value = "Norma (Sistem əhəmiyyətli)"
$ws2.Range("C1").Value = "Norma (Banklar istisna)"
$ws2.Range("D1").Value = "Fakt"

# The remaining values look like percentages ("6.0%") or contain a "%"
# sign ("minimum 5%"); force them to be stored as plain text instead of
# letting Excel auto-convert them to numeric percentages, then strip the
# temporary text number-format back off so the cells keep the workbook's
# default (unstyled) cell style.
$percentCells = $ws2.Range("B2:D4")
$percentCells.NumberFormat = "@"

$ws2.Range("B2").Value = "6.0%"
$ws2.Range("C2").Value = "5.0%"
$ws2.Range("D2").Value = "17.0%"

$ws2.Range("B3").Value = "11.0%"
$ws2.Range("C3").Value = "9.0%"
$ws2.Range("D3").Value = "21.0%"

$ws2.Range("B4").Value = "minimum 5%"
$ws2.Range("C4").Value = "minimum 4%"
$ws2.Range("D4").Value = "10.0%"

$percentCells.ClearFormats()

# Match the header-row styling (bold font, thin borders, centered/top
# aligned) used for the matching header cells on Table_1 by copying the
# existing formatting so the same style record is reused.
$srcHeader = $ws1.Range("A1:B1")
$srcHeader.Copy()
$dstHeader = $ws2.Range("A1:D1")
$dstHeader.PasteSpecial(-4122)

$ws2.Range("A1").Select()
